# edit.ps1 - Applies "Added Domains to Explination" changes to ERExplination.docx
#
# Strategy: locate each target paragraph by its (pre-edit) index, and replace
# the paragraph's Range content with freshly authored OOXML via Range.InsertXML.
# This lets us split/merge runs exactly as required without fighting Find/Replace
# run-boundary semantics. We process paragraphs from the bottom of the document
# upward so that paragraph-count-changing edits (the single new paragraph near
# the top) never invalidate indices we still need to use.

$d = $word.ActiveDocument
$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-ParaXml($doc, $index, $innerXml) {
    $p = $doc.Paragraphs.Item($index)
    $r = $p.Range
    $xml = "<w:p $W>" + $innerXml + "</w:p>"
    $r.InsertXML($xml)
}

# --- 1. "Cardinality:" (Works relationship block, paragraph 30) gains the
#        lastRenderedPageBreak marker that used to sit on the next paragraph.
$cardinalityPPr = "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr>"
$cardinalityRuns = "<w:r><w:lastRenderedPageBreak/><w:t>Cardinality:</w:t></w:r>"
$inner1 = $cardinalityPPr + $cardinalityRuns
Set-ParaXml $d 30 $inner1

# --- 2. "User (1,1): Each workout is specific..." (paragraph 31) loses the
#        lastRenderedPageBreak marker.
$userPPr = "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr>"
$userRuns = "<w:r><w:t>User (1,</w:t></w:r>" +
            "<w:r><w:t>1</w:t></w:r>" +
            "<w:r><w:t xml:space=`"preserve`">): </w:t></w:r>" +
            "<w:r><w:t>Each workout is specific to the user that accomplishes them.</w:t></w:r>"
$inner2 = $userPPr + $userRuns
Set-ParaXml $d 31 $inner2

# --- 3. "Recommended Daily Dose: Self explanatory" (paragraph 23)
$lvl1PPr = "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>"
$infinityRpr = "<w:rPr><w:rFonts w:cstheme=`"minorHAnsi`"/><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"FFFFFF`"/></w:rPr>"
$rddRuns = "<w:r><w:t>Recommended Daily Dose</w:t></w:r>" +
           "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
           "<w:r><w:t>(</w:t></w:r>" +
           "<w:r><w:t>1</w:t></w:r>" +
           "<w:r><w:t>-</w:t></w:r>" +
           "<w:r>$infinityRpr<w:t>&#8734;</w:t></w:r>" +
           "<w:r><w:t>)</w:t></w:r>" +
           "<w:r><w:t>: Self explanatory</w:t></w:r>"
$inner3 = $lvl1PPr + $rddRuns
Set-ParaXml $d 23 $inner3

# --- 4. "Calories: the calories per serving" (paragraph 20)
$caloriesRuns = "<w:r><w:t>Calories</w:t></w:r>" +
                "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
                "<w:r><w:t>(</w:t></w:r>" +
                "<w:r><w:t>1</w:t></w:r>" +
                "<w:r><w:t>-</w:t></w:r>" +
                "<w:r>$infinityRpr<w:t>&#8734;</w:t></w:r>" +
                "<w:r><w:t>)</w:t></w:r>" +
                "<w:r><w:t xml:space=`"preserve`">: the calories </w:t></w:r>" +
                "<w:r><w:t>ounce</w:t></w:r>"
$inner4 = $lvl1PPr + $caloriesRuns
Set-ParaXml $d 20 $inner4

# --- 5. "Grams Per Serving: Self explanatory" (paragraph 19)
$gpsRuns = "<w:r><w:t>Grams Per Serving</w:t></w:r>" +
           "<w:r><w:t xml:space=`"preserve`"> (1-</w:t></w:r>" +
           "<w:r>$infinityRpr<w:t>&#8734;</w:t></w:r>" +
           "<w:r><w:t>)</w:t></w:r>" +
           "<w:r><w:t>: Self explanatory</w:t></w:r>"
$inner5 = $lvl1PPr + $gpsRuns
Set-ParaXml $d 19 $inner5

# --- 6. "ID: The unique ID of the workout" (paragraph 10)
$idWorkoutRuns = "<w:r><w:t>ID</w:t></w:r>" +
                 "<w:r><w:t xml:space=`"preserve`"> (0-</w:t></w:r>" +
                 "<w:r>$infinityRpr<w:t>&#8734;</w:t></w:r>" +
                 "<w:r><w:t>)</w:t></w:r>" +
                 "<w:r><w:t>: The unique ID of the workout</w:t></w:r>"
$inner6 = $lvl1PPr + $idWorkoutRuns
Set-ParaXml $d 10 $inner6

# --- 7. "Weight: The weight of the weight entity" (paragraph 7)
$weightRuns = "<w:r><w:t>Weight</w:t></w:r>" +
              "<w:r><w:t xml:space=`"preserve`"> (0-</w:t></w:r>" +
              "<w:r>$infinityRpr<w:t>&#8734;</w:t></w:r>" +
              "<w:r><w:t>)</w:t></w:r>" +
              "<w:r><w:t>: The weight of the weight entity</w:t></w:r>"
$inner7 = $lvl1PPr + $weightRuns
Set-ParaXml $d 7 $inner7

# --- 8. "ID: This is the identifier for the User entity. Given that..." (paragraph 3)
$idUserRuns = "<w:r><w:t>ID</w:t></w:r>" +
              "<w:r><w:t xml:space=`"preserve`"> (0-</w:t></w:r>" +
              "<w:r>$infinityRpr<w:t>&#8734;</w:t></w:r>" +
              "<w:r><w:t>)</w:t></w:r>" +
              "<w:r><w:t xml:space=`"preserve`">: </w:t></w:r>" +
              "<w:r><w:t>This is the identifier for the User entity</w:t></w:r>" +
              "<w:r><w:t>. Given that many users could have the same name a unique ID will be generated for each one.</w:t></w:r>"
$inner8 = $lvl1PPr + $idUserRuns
Set-ParaXml $d 3 $inner8

# --- 9. New italic "Domain of attributes are in ()" paragraph, inserted right
#        after the "ENTITIES & ATTRIBUTES" heading (paragraph 1).
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$domainPara = $d.Paragraphs.Item(2)
$domainPPr = "<w:pPr><w:ind w:left=`"720`" w:hanging=`"360`"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr>"
$italicRpr = "<w:rPr><w:i/><w:iCs/></w:rPr>"
$domainRuns = "<w:r>$italicRpr<w:t>Domain</w:t></w:r>" +
              "<w:r>$italicRpr<w:t xml:space=`"preserve`"> of attributes</w:t></w:r>" +
              "<w:r>$italicRpr<w:t xml:space=`"preserve`"> </w:t></w:r>" +
              "<w:r>$italicRpr<w:t>are</w:t></w:r>" +
              "<w:r>$italicRpr<w:t xml:space=`"preserve`"> in ()</w:t></w:r>"
$domainXml = "<w:p $W>" + $domainPPr + $domainRuns + "</w:p>"
$domainPara.Range.InsertXML($domainXml)

Write-Host "Done."
